$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2204.2856
$ws.Range("I34").Value = 2204.2856
$ws.Range("K34").Value = 2204.2856
$ws.Range("M34").Value = -2001.2856
$ws.Range("H36").Value = 2204.2856
$ws.Range("I36").Value = 2204.2856
$ws.Range("K36").Value = 2204.2856
$ws.Range("M36").Value = -1489.2856
$ws.Range("H62").Value = 26519.32
$ws.Range("I62").Value = 6351.9414
$ws.Range("K62").Value = 6351.9414
$ws.Range("M62").Value = -5727.9414
$ws.Range("H64").Value = 4140.4165
$ws.Range("J64").Value = 4162.273
$ws.Range("L64").Value = 4162.273
$ws.Range("N64").Value = -4658.273
$ws.Range("H65").Value = 26519.32
$ws.Range("I65").Value = 6351.9414
$ws.Range("K65").Value = 31759.707
$ws.Range("M65").Value = -28639.707
$ws.Range("H67").Value = 4140.4165
$ws.Range("J67").Value = 4162.273
$ws.Range("L67").Value = 4162.273
$ws.Range("N67").Value = -5878.273
$ws.Range("H86").Value = 15466962
$ws.Range("J86").Value = 20106100
$ws.Range("L86").Value = 20106100
$ws.Range("N86").Value = -20108346
$ws.Range("H89").Value = 15466962
$ws.Range("J89").Value = 20106100
$ws.Range("L89").Value = 100530500
$ws.Range("N89").Value = -100541732
$ws.Range("H96").Value = 2001.625
$ws.Range("I96").Value = 3085.25
$ws.Range("K96").Value = 9255.75
$ws.Range("M96").Value = -7882.75
$ws.Range("H98").Value = 1687.0588
$ws.Range("I98").Value = 855.6
$ws.Range("K98").Value = 855.6
$ws.Range("M98").Value = 642.4
$ws.Range("H113").Value = 90913160
$ws.Range("J113").Value = 4199.778
$ws.Range("L113").Value = 4199.778
$ws.Range("N113").Value = -10707.778
$ws.Range("H116").Value = 53136216
$ws.Range("I116").Value = 27906112
$ws.Range("K116").Value = 27906112
$ws.Range("M116").Value = -27902670
$ws.Range("H122").Value = 1687.0588
$ws.Range("I122").Value = 855.6
$ws.Range("K122").Value = 2566.8
$ws.Range("M122").Value = -116.8000000000002
$ws.Range("H127").Value = 1701
$ws.Range("I127").Value = 495
$ws.Range("J127").Value = 2183.4
$ws.Range("K127").Value = 1485
$ws.Range("L127").Value = 6550.200000000001
$ws.Range("M127").Value = 3475
$ws.Range("N127").Value = -16470.2
$ws.Range("H135").Value = 843.6667
$ws.Range("I135").Value = 786.625
$ws.Range("K135").Value = 7079.625
$ws.Range("M135").Value = -4544.625
$ws.Range("H137").Value = 4242.1035
$ws.Range("I137").Value = 3232.7368
$ws.Range("J137").Value = 6159.9
$ws.Range("K137").Value = 9698.2104
$ws.Range("L137").Value = 18479.7
$ws.Range("M137").Value = -7148.2104
$ws.Range("N137").Value = -23579.7
$ws.Range("H138").Value = 3580.5417
$ws.Range("J138").Value = 3572.4167
$ws.Range("L138").Value = 10717.2501
$ws.Range("N138").Value = -20997.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1639.4231
$ws.Range("I2").Value = 1042.2858
$ws.Range("J2").Value = 4147.4
$ws.Range("K2").Value = 1042.2858
$ws.Range("L2").Value = 4147.4
$ws.Range("M2").Value = -929.2858000000001
$ws.Range("N2").Value = -4373.4
$ws.Range("H32").Value = 16265.525
$ws.Range("I32").Value = 8938.574000000001
$ws.Range("J32").Value = 39788.895
$ws.Range("K32").Value = 8938.574000000001
$ws.Range("L32").Value = 39788.895
$ws.Range("M32").Value = -8651.574000000001
$ws.Range("N32").Value = -40362.895
$ws.Range("H61").Value = 3306.652
$ws.Range("I61").Value = 2377.65
$ws.Range("K61").Value = 2377.65
$ws.Range("M61").Value = -2165.65
$ws.Range("H74").Value = 40004756
$ws.Range("J74").Value = 4932.7334
$ws.Range("L74").Value = 4932.7334
$ws.Range("N74").Value = -6680.7334
$ws.Range("H77").Value = 40004756
$ws.Range("J77").Value = 4932.7334
$ws.Range("L77").Value = 24663.667
$ws.Range("N77").Value = -33399.667
$ws.Range("H116").Value = 1639.4231
$ws.Range("I116").Value = 1042.2858
$ws.Range("J116").Value = 4147.4
$ws.Range("K116").Value = 1042.2858
$ws.Range("L116").Value = 4147.4
$ws.Range("M116").Value = 1251.7142
$ws.Range("N116").Value = -8735.4
$ws.Range("H132").Value = 2996.8572
$ws.Range("I132").Value = 2551.7742
$ws.Range("K132").Value = 7655.3226
$ws.Range("M132").Value = -5125.3226
$ws.Range("H136").Value = 3306.652
$ws.Range("I136").Value = 2377.65
$ws.Range("K136").Value = 7132.950000000001
$ws.Range("M136").Value = -4582.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1639.4231
$ws.Range("I3").Value = 1042.2858
$ws.Range("J3").Value = 4147.4
$ws.Range("K3").Value = 1042.2858
$ws.Range("L3").Value = 4147.4
$ws.Range("M3").Value = -928.2858000000001
$ws.Range("N3").Value = -4375.4
$ws.Range("H99").Value = 3448.2354
$ws.Range("I99").Value = 3074.6667
$ws.Range("K99").Value = 3074.6667
$ws.Range("M99").Value = -1576.6667
$ws.Range("H107").Value = 62502576
$ws.Range("I107").Value = 2836.6667
$ws.Range("J107").Value = 100002424
$ws.Range("K107").Value = 2836.6667
$ws.Range("L107").Value = 100002424
$ws.Range("M107").Value = -916.6667000000002
$ws.Range("N107").Value = -100006264

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4268.1924
$ws.Range("I31").Value = 3737.25
$ws.Range("K31").Value = 3737.25
$ws.Range("M31").Value = -3442.25
$ws.Range("H34").Value = 4268.1924
$ws.Range("I34").Value = 3737.25
$ws.Range("K34").Value = 3737.25
$ws.Range("M34").Value = -3535.25
$ws.Range("H132").Value = 201640.34
$ws.Range("I132").Value = 1223.5526
$ws.Range("K132").Value = 3670.6578
$ws.Range("M132").Value = -1140.6578

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 572053.5600000001
$ws.Range("I68").Value = 1333666.6
$ws.Range("J68").Value = 843.75
$ws.Range("K68").Value = 4000999.8
$ws.Range("L68").Value = 2531.25
$ws.Range("M68").Value = -4000188.8
$ws.Range("N68").Value = -4153.25
$ws.Range("H71").Value = 572053.5600000001
$ws.Range("I71").Value = 1333666.6
$ws.Range("J71").Value = 843.75
$ws.Range("K71").Value = 12002999.4
$ws.Range("L71").Value = 7593.75
$ws.Range("M71").Value = -11998943.4
$ws.Range("N71").Value = -15705.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3447.7144
$ws.Range("I80").Value = 3199
$ws.Range("J80").Value = 3547.2
$ws.Range("K80").Value = 3199
$ws.Range("L80").Value = 3547.2
$ws.Range("M80").Value = -2201
$ws.Range("N80").Value = -5543.2
$ws.Range("H83").Value = 3447.7144
$ws.Range("I83").Value = 3199
$ws.Range("J83").Value = 3547.2
$ws.Range("K83").Value = 15995
$ws.Range("L83").Value = 17736
$ws.Range("M83").Value = -11003
$ws.Range("N83").Value = -27720
$ws.Range("H102").Value = 2164.2896
$ws.Range("I102").Value = 1045.6522
$ws.Range("J102").Value = 3879.5334
$ws.Range("K102").Value = 1045.6522
$ws.Range("L102").Value = 3879.5334
$ws.Range("M102").Value = 576.3478
$ws.Range("N102").Value = -7123.5334
$ws.Range("H132").Value = 4377.6924
$ws.Range("I132").Value = 4628.1816
$ws.Range("K132").Value = 13884.5448
$ws.Range("M132").Value = -11354.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3169.3044
$ws.Range("I40").Value = 2778.4285
$ws.Range("K40").Value = 2778.4285
$ws.Range("M40").Value = -2642.4285
$ws.Range("H109").Value = 105000.5
$ws.Range("J109").Value = 105000.5
$ws.Range("L109").Value = 105000.5
$ws.Range("N109").Value = -107774.5
$ws.Range("H132").Value = 16954096
$ws.Range("J132").Value = 8211
$ws.Range("L132").Value = 24633
$ws.Range("N132").Value = -29693
$ws.Range("H133").Value = 99999.5
$ws.Range("J133").Value = 99999.5
$ws.Range("L133").Value = 99999.5
$ws.Range("N133").Value = -105059.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1950.5
$ws.Range("I122").Value = 1647.6923
$ws.Range("K122").Value = 4943.0769
$ws.Range("M122").Value = -2493.0769
$ws.Range("H132").Value = 512188.7
$ws.Range("I132").Value = 774399.75
$ws.Range("J132").Value = 25225.285
$ws.Range("K132").Value = 2323199.25
$ws.Range("L132").Value = 75675.855
$ws.Range("M132").Value = -2320669.25
$ws.Range("N132").Value = -80735.855
$ws.Range("H136").Value = 8581.264999999999
$ws.Range("I136").Value = 11645.3125
$ws.Range("J136").Value = 5857.6665
$ws.Range("K136").Value = 34935.9375
$ws.Range("L136").Value = 17572.9995
$ws.Range("M136").Value = -32385.9375
$ws.Range("N136").Value = -22672.9995
